$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q3").Value = "'30491.43"
$ws.Range("Q4").Value = "'1140.99"
$ws.Range("Q5").Value = "'17412.80"
$ws.Range("Q6").Value = "'65114.86"
$ws.Range("Q7").Value = "'2003.07"
$ws.Range("Q8").Value = "'1892.25"
$ws.Range("Q9").Value = "'1386.30"
$ws.Range("Q10").Value = "'34680.21"
$ws.Range("Q11").Value = "'7259.67"
$ws.Range("Q12").Value = "'142.84"
$ws.Range("Q13").Value = "'155.15"
$ws.Range("Q14").Value = "'28.67"
$ws.Range("Q15").Value = "'1972.19"
$ws.Range("Q16").Value = "'9977.52"
$ws.Range("Q17").Value = "'236328.47"
$ws.Range("Q18").Value = "'80412.84"
$ws.Range("Q19").Value = "'20959.92"
$ws.Range("Q20").Value = "'20.55"
$ws.Range("Q21").Value = "'884699.58"
$ws.Range("Q22").Value = "'11.09"
$ws.Range("Q23").Value = "'99251.99"
